$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark from its original location (right
#    before "Informatyka, Wydział Informatyki PUT").
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 2. Append " .NET" (Strong style, not bold) as its own run to the
#    paragraph "Omówienie stosowanych technologii".
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Omówienie stosowanych technologii`r") {
        $r = $p.Range
        $r.MoveEnd(1, -1)
        $splitPos = $r.End

        $endPt = $r.Duplicate
        $endPt.Collapse(0)
        $endPt.InsertAfter(" .NET")

        # Force the newly appended text into its own run by briefly
        # bookmarking the split point (now interior to the paragraph,
        # so this is safe) and removing the bookmark again.
        $splitRange = $d.Range($splitPos, $splitPos)
        $d.Bookmarks.Add("TMP_SPLIT_1", $splitRange)
        $d.Bookmarks("TMP_SPLIT_1").Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 3. Delete the (old) "Izolacja modułów" paragraph that currently sits
#    right after "Aktualizacje" -- its text is being relocated to the
#    paragraph that used to read "Rozproszone repozytorium kodu".
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Izolacja modułów`r") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 4. Rename "Rozproszone repozytorium kodu" -> "Izolacja modułów".
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Rozproszone repozytorium kodu`r") {
        $r = $p.Range
        $r.MoveEnd(1, -1)
        $r.Text = "Izolacja modułów"
        break
    }
}

# ---------------------------------------------------------------------
# 5. The old "Literatura" bullet becomes "Omówienie" + bookmark +
#    " zarządzania projektem", and a brand new "Literatura" bullet is
#    inserted right after it.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Literatura`r") {
        $r = $p.Range
        $r.MoveEnd(1, -1)
        $r.Text = "Omówienie"

        $bmPos = $r.Start + 9   # length of "Omówienie"

        $endPt = $r.Duplicate
        $endPt.Collapse(0)
        $endPt.InsertAfter(" zarządzania projektem")

        $bmRange = $d.Range($bmPos, $bmPos)
        $d.Bookmarks.Add("_GoBack", $bmRange)

        $newPara = $p.Range.InsertParagraphAfter()
        $newRange = $d.Paragraphs.Item($p.Index + 1).Range
        $newRange.MoveEnd(1, -1)
        $newRange.Text = "Literatura"
        break
    }
}
